# integrate web testing with this framework
#
# Applies the changes described by the commit:
#  - MAIN_CONTROLLER (sheet1): drop the ApplicationType column (F), fix up
#    row 4 (now a passing "Ishine" run) and append a new row 5 (FOS run),
#    and flag the failing row 3 with a red fill.
#  - DATASHEET (sheet2): point the existing "FOS" row at FOS8.xlsx and add a
#    new data row for the Ishine application/spreadsheet.
#  - Selection / active-tab bookkeeping so DATASHEET ends up the visible tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws2 = $wb.Worksheets.Item("DATASHEET")
$ws4 = $wb.Worksheets.Item("MAIL_SEND")

# ---------------------------------------------------------------------
# MAIN_CONTROLLER: remove the ApplicationType column (F) entirely
# ---------------------------------------------------------------------
$ws1.Range("F1:F4").Delete(-4159)

# Row 4: this run now succeeds (Y) and targets the new "Ishine" app
$ws1.Range("B4").Value = "Y"
$ws1.Range("D4").Value = "Ishine"
$ws1.Range("E4").Value = "Ishine"

# New row 5: an additional FOS run
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Y"
$ws1.Range("C5").Value = "local"
$ws1.Range("D5").Value = "FOS"
$ws1.Range("E5").Value = "FOS"

# Highlight the failing row (row 3, RunStatus = N) with a red fill
$ws1.Range("B3").Interior.Color = 255

# ---------------------------------------------------------------------
# DATASHEET: add the Ishine data row first (so "Ishine.xlsx" gets
# registered in the shared-string table before "FOS8.xlsx"), then update
# the existing FOS row to point at FOS8.xlsx
# ---------------------------------------------------------------------
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "N"
$ws2.Range("C4").Value = "Ishine"
$ws2.Range("D4").Value = "Ishine.xlsx"
$ws2.Range("E4").Value = 2
$ws2.Range("F4").Value = 20

$ws2.Range("D3").Value = "FOS8.xlsx"

# ---------------------------------------------------------------------
# Selections
# ---------------------------------------------------------------------
$ws4.Range("A3").EntireRow.Select()

$ws1.Range("I7").Select()

# DATASHEET becomes the active / selected tab
$ws2.Range("E11").Select()
